$wb = $excel.ActiveWorkbook

# Row 28 on sheet ALC (diff @ @@ -2010,25 +2010,25 @@)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 688.7059
$ws.Range("J28").Value = 837.5
$ws.Range("L28").Value = 837.5
$ws.Range("N28").Value = -1807.5

# Row 70 on sheet ALC (diff @ @@ -4056,22 +4056,25 @@)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1433
$ws.Range("I70").Value = 300
$ws.Range("J70").Value = 1999.5
$ws.Range("K70").Value = 900
$ws.Range("L70").Value = 5998.5
$ws.Range("M70").Value = -630
$ws.Range("N70").Value = -6538.5

# Row 73 on sheet ALC (diff @ @@ -4200,22 +4203,25 @@)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 1433
$ws.Range("I73").Value = 300
$ws.Range("J73").Value = 1999.5
$ws.Range("K73").Value = 900
$ws.Range("L73").Value = 5998.5
$ws.Range("M73").Value = 36
$ws.Range("N73").Value = -7870.5

# Row 112 on sheet ALC (diff @ @@ -6129,22 +6135,22 @@)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 2465.8572
$ws.Range("I112").Value = 1361.6666
$ws.Range("K112").Value = 4084.9998
$ws.Range("M112").Value = -2976.9998

# Row 132 on sheet ALC (diff @ @@ -7103,22 +7109,22 @@)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 5190.3
$ws.Range("I132").Value = 3129.4285
$ws.Range("K132").Value = 9388.2855
$ws.Range("M132").Value = -6858.2855

# Row 138 on sheet ALC (diff @ @@ -7403,22 +7409,22 @@)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 8299.950000000001
$ws.Range("I138").Value = 2249.75
$ws.Range("K138").Value = 6749.25
$ws.Range("M138").Value = -1609.25

# Row 4 on sheet ARM (diff @ @@ -7794,25 +7800,25 @@)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 82
$ws.Range("I4").Value = 89
$ws.Range("J4").Value = 75
$ws.Range("K4").Value = 89
$ws.Range("L4").Value = 75
$ws.Range("M4").Value = 27
$ws.Range("N4").Value = -307

# Row 6 on sheet ARM (diff @ @@ -7898,22 +7904,19 @@)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("N6").ClearContents()

# Row 74 on sheet ARM (diff @ @@ -11167,22 +11170,22 @@)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1103.6154
$ws.Range("I74").Value = 949.7273
$ws.Range("K74").Value = 949.7273
$ws.Range("M74").Value = -75.72730000000001

# Row 77 on sheet ARM (diff @ @@ -11311,22 +11314,22 @@)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1103.6154
$ws.Range("I77").Value = 949.7273
$ws.Range("K77").Value = 4748.636500000001
$ws.Range("M77").Value = -380.6365000000005

# Row 110 on sheet ARM (diff @ @@ -12883,22 +12886,19 @@)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 0
$ws.Range("I110").Value = 0
$ws.Range("K110").Value = 0
$ws.Range("M110").ClearContents()

# Row 122 on sheet ARM (diff @ @@ -13450,22 +13450,22 @@)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2000.3334
$ws.Range("I122").Value = 2000.3334
$ws.Range("K122").Value = 6001.0002
$ws.Range("M122").Value = -3551.0002

# Row 107 on sheet BSM (diff @ @@ -19525,22 +19525,22 @@)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 196.5
$ws.Range("I107").Value = 196.5
$ws.Range("K107").Value = 196.5
$ws.Range("M107").Value = 1723.5

# Row 2 on sheet CRP (diff @ @@ -21259,25 +21259,25 @@)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 6576.25
$ws.Range("I2").Value = 1000
$ws.Range("J2").Value = 8435
$ws.Range("K2").Value = 1000
$ws.Range("L2").Value = 8435
$ws.Range("M2").Value = -887
$ws.Range("N2").Value = -8661

# Row 7 on sheet CRP (diff @ @@ -21504,22 +21504,22 @@)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 1677.8
$ws.Range("I7").Value = 997
$ws.Range("K7").Value = 997
$ws.Range("M7").Value = -884

# Row 22 on sheet CRP (diff @ @@ -22245,22 +22245,25 @@)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 484.14285
$ws.Range("I22").Value = 445
$ws.Range("J22").Value = 499.8
$ws.Range("K22").Value = 445
$ws.Range("L22").Value = 499.8
$ws.Range("M22").Value = -95
$ws.Range("N22").Value = -1199.8

# Row 29 on sheet CRP (diff @ @@ -22594,22 +22597,22 @@)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H29").Value = 2000
$ws.Range("J29").Value = 2000
$ws.Range("L29").Value = 2000
$ws.Range("N29").Value = -2586

# Row 31 on sheet CRP (diff @ @@ -22692,25 +22695,25 @@)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6298.778
$ws.Range("I31").Value = 6138.8
$ws.Range("J31").Value = 6498.75
$ws.Range("K31").Value = 6138.8
$ws.Range("L31").Value = 6498.75
$ws.Range("M31").Value = -5843.8
$ws.Range("N31").Value = -7088.75

# Row 34 on sheet CRP (diff @ @@ -22845,25 +22848,25 @@)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 6298.778
$ws.Range("I34").Value = 6138.8
$ws.Range("J34").Value = 6498.75
$ws.Range("K34").Value = 6138.8
$ws.Range("L34").Value = 6498.75
$ws.Range("M34").Value = -5936.8
$ws.Range("N34").Value = -6902.75

# Row 35 on sheet CRP (diff @ @@ -22897,22 +22900,25 @@)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H35").Value = 2644.9092
$ws.Range("I35").Value = 2529.4
$ws.Range("J35").Value = 3800
$ws.Range("K35").Value = 2529.4
$ws.Range("L35").Value = 3800
$ws.Range("M35").Value = -2235.4
$ws.Range("N35").Value = -4388

# Row 58 on sheet CRP (diff @ @@ -23994,22 +24000,25 @@)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1749.875
$ws.Range("I58").Value = 1857
$ws.Range("J58").Value = 1000
$ws.Range("K58").Value = 1857
$ws.Range("L58").Value = 1000
$ws.Range("M58").Value = -1654
$ws.Range("N58").Value = -1406

# Row 99 on sheet CRP (diff @ @@ -25979,25 +25988,25 @@)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2846.6
$ws.Range("J99").Value = 3372.25
$ws.Range("L99").Value = 3372.25
$ws.Range("N99").Value = -6368.25

# Row 105 on sheet CRP (diff @ @@ -26267,22 +26276,22 @@)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 4254.385
$ws.Range("I105").Value = 2851.125
$ws.Range("K105").Value = 2851.125
$ws.Range("M105").Value = -1104.125

# Row 107 on sheet CRP (diff @ @@ -26368,22 +26377,22 @@)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 921
$ws.Range("I107").Value = 756.5
$ws.Range("K107").Value = 756.5
$ws.Range("M107").Value = 1163.5

# Row 122 on sheet CRP (diff @ @@ -27073,25 +27082,25 @@)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1754.125
$ws.Range("I122").Value = 600
$ws.Range("J122").Value = 2138.8333
$ws.Range("K122").Value = 1800
$ws.Range("L122").Value = 6416.499899999999
$ws.Range("M122").Value = 650
$ws.Range("N122").Value = -11316.4999

# Row 126 on sheet CRP (diff @ @@ -27266,25 +27275,25 @@)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 2846.6
$ws.Range("J126").Value = 3372.25
$ws.Range("L126").Value = 10116.75
$ws.Range("N126").Value = -15056.75

# Row 132 on sheet CRP (diff @ @@ -27551,22 +27560,19 @@)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

# Row 134 on sheet CRP (diff @ @@ -27646,22 +27652,22 @@)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1987.4286
$ws.Range("I134").Value = 1983
$ws.Range("K134").Value = 5949
$ws.Range("M134").Value = -3414

# Row 136 on sheet CRP (diff @ @@ -27747,22 +27753,25 @@)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1749.875
$ws.Range("I136").Value = 1857
$ws.Range("J136").Value = 1000
$ws.Range("K136").Value = 5571
$ws.Range("L136").Value = 3000
$ws.Range("M136").Value = -3021
$ws.Range("N136").Value = -8100

# Row 104 on sheet CUL (diff @ @@ -33169,7 +33178,7 @@)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H104").Value = 2200

# Row 70 on sheet GSM (diff @ @@ -38424,25 +38433,19 @@)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("M70").ClearContents()
$ws.Range("N70").ClearContents()

# Row 73 on sheet GSM (diff @ @@ -38568,25 +38571,19 @@)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("M73").ClearContents()
$ws.Range("N73").ClearContents()

# Row 102 on sheet GSM (diff @ @@ -39959,22 +39956,22 @@)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1015.3571
$ws.Range("I102").Value = 1015.3571
$ws.Range("K102").Value = 1015.3571
$ws.Range("M102").Value = 606.6429000000001

# Row 113 on sheet GSM (diff @ @@ -40474,23 +40471,26 @@)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 4485
$ws.Range("J113").Value = 4999
$ws.Range("L113").Value = 4999
$ws.Range("N113").Value = -9339

# Row 122 on sheet GSM (diff @ @@ -40900,22 +40900,22 @@)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 5451.1665
$ws.Range("I122").Value = 5451.1665
$ws.Range("K122").Value = 16353.4995
$ws.Range("M122").Value = -13903.4995

# Row 123 on sheet GSM (diff @ @@ -40949,22 +40949,22 @@)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 33000
$ws.Range("J123").Value = 33000
$ws.Range("L123").Value = 33000
$ws.Range("N123").Value = -37900

# Row 126 on sheet GSM (diff @ @@ -41093,25 +41093,22 @@)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 5400
$ws.Range("I126").Value = 5400
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 16200
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -13730
$ws.Range("N126").ClearContents()

# Row 132 on sheet GSM (diff @ @@ -41375,22 +41372,22 @@)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3582.4666
$ws.Range("I132").Value = 3499.1035
$ws.Range("K132").Value = 10497.3105
$ws.Range("M132").Value = -7967.3105

# Row 22 on sheet LTW (diff @ @@ -42918,25 +42915,25 @@)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2500
$ws.Range("I22").Value = 1000
$ws.Range("J22").Value = 3000
$ws.Range("K22").Value = 1000
$ws.Range("L22").Value = 3000
$ws.Range("M22").Value = -705
$ws.Range("N22").Value = -3590

# Row 27 on sheet LTW (diff @ @@ -43157,25 +43154,25 @@)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 2500
$ws.Range("I27").Value = 1000
$ws.Range("J27").Value = 3000
$ws.Range("K27").Value = 1000
$ws.Range("L27").Value = 3000
$ws.Range("M27").Value = -893
$ws.Range("N27").Value = -3214

# Row 93 on sheet LTW (diff @ @@ -46325,19 +46322,25 @@)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 492
$ws.Range("I93").Value = 600
$ws.Range("J93").Value = 384
$ws.Range("K93").Value = 600
$ws.Range("L93").Value = 384
$ws.Range("M93").Value = 648
$ws.Range("N93").Value = -2880

# Row 100 on sheet LTW (diff @ @@ -46650,23 +46653,26 @@)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 1424.25
$ws.Range("J100").Value = 1398
$ws.Range("L100").Value = 1398
$ws.Range("N100").Value = -2480

# Row 132 on sheet LTW (diff @ @@ -48182,22 +48188,25 @@)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 7601.2
$ws.Range("I132").Value = 3890
$ws.Range("J132").Value = 8529
$ws.Range("K132").Value = 11670
$ws.Range("L132").Value = 25587
$ws.Range("M132").Value = -9140
$ws.Range("N132").Value = -30647

# Row 122 on sheet WVR (diff @ @@ -54532,25 +54541,25 @@)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3399.7
$ws.Range("I122").Value = 3124.75
$ws.Range("J122").Value = 4499.5
$ws.Range("K122").Value = 9374.25
$ws.Range("L122").Value = 13498.5
$ws.Range("M122").Value = -6924.25
$ws.Range("N122").Value = -18398.5
